# v1.1 verify that the previous comments modified
# close publish article wireframe review and verify the updates

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LH_WF_PUBLISHARTICLE_REVIEW")
$ws2 = $wb.Worksheets.Item("VERSION-HISTORY")

# --- Close out the review: mark all three reviewer-verification rows as "closed" ---
$ws1.Range("I2").Value = "closed"
$ws1.Range("I3").Value = "closed"
$ws1.Range("I4").Value = "closed"

# Row 3's "Owner Status" cell (H3) had its highlight fill cleared once the
# review item was closed out.
$ws1.Range("H3").Interior.Pattern = -4142

# --- Record the new version history entry for this review closure ---
$ws2.Range("A3").Value = "v1.1"
$ws2.Range("B3").Value = "Ahmed Abuzaid"
$ws2.Range("C3").Value = "close the review status "
$ws2.Range("D3").Value = 45776

# --- Restore selections to match where the reviewer finished working ---
$ws2.Range("C3").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("I4").Select() | Out-Null
